# Apply cryptos list update (Mon Apr  1 20:53:41 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.661.70"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "3.485.00"
$ws.Range("E3").Value = "  -4.10%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.74"
$ws.Range("E5").Value = "  -4.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.87"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("D8").Value = "3.473.26"
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -7.57%  "
$ws.Range("E11").Value = "  -4.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.48"
$ws.Range("E12").Value = "  -4.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("E13").Value = "  -6.38%  "
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").Value = "4.036.44"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "648.60"
$ws.Range("E16").Value = "  -4.31%  "
$ws.Range("D17").Value = "69.474.94"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "3.480.80"
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("E19").Value = "  -5.44%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  -4.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.948"
$ws.Range("E22").Value = "  -5.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.16"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.24"
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.87"
$ws.Range("E25").Value = "  -6.45%  "
$ws.Range("E26").Value = "  -7.42%  "
$ws.Range("E27").Value = "  -4.28%  "
$ws.Range("E28").Value = "  -4.07%  "
$ws.Range("E29").Value = "  -4.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.61"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.27"
$ws.Range("E31").Value = "  -8.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.74"
$ws.Range("E32").Value = "  -5.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.62"
$ws.Range("E33").Value = "  -4.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("E34").Value = "  -5.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "61.03"
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("D36").Value = "3.730.03"
$ws.Range("E36").Value = "  -5.95%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "527.01"
$ws.Range("E37").Value = "  +4.13%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -8.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.54"
$ws.Range("E43").Value = "  +68.08%  "
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.33"
$ws.Range("E45").Value = "  -6.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0443"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("E48").Value = "  -8.27%  "
$ws.Range("E49").Value = "  -4.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -5.96%  "
